$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - update existing "Nunawading" entry to "Abbotsford" entry, mark as old
$ws.Cells.Item(2,1).Value = "Abbotsford"
$ws.Cells.Item(2,2).Value = "Bodriggy Brewing Company  245 Johnston Street, Abbotsford VIC 3067"
$ws.Cells.Item(2,3).Value = "28/12/20 2:50pm-5:30pm"
$ws.Cells.Item(2,4).Value = "Case dined at venue"
$ws.Cells.Item(2,5).Value = "old"

# Row 3 - Albert Park
$ws.Cells.Item(3,1).Value = "Albert Park"
$ws.Cells.Item(3,2).Value = "The Guilty Moose Cafe  143 Victoria Avenue, Albert Park VIC 3206"
$ws.Cells.Item(3,3).Value = "21/12/20 1pm-1:30pm"
$ws.Cells.Item(3,4).Value = "Case ate at cafe"
$ws.Cells.Item(3,5).Value = "old"

# Row 4 - Hampton
$ws.Cells.Item(4,1).Value = "Hampton"
$ws.Cells.Item(4,2).Value = "Merrymen Cafe, 2 Small Street, Hampton VIC"
$ws.Cells.Item(4,3).Value = "28-12-2020 1:30pm-2:30pm"
$ws.Cells.Item(4,4).Value = "Case ate in store"
$ws.Cells.Item(4,5).Value = "old"

# Row 5 - Springvale (IKEA)
$ws.Cells.Item(5,1).Value = "Springvale"
$ws.Cells.Item(5,2).Value = "IKEA Springvale, 917 Princes Hwy"
$ws.Cells.Item(5,3).Value = "30/12/20, 4:00pm-6:30pm"
$ws.Cells.Item(5,4).Value = "Case shopped at store and dined at cafe"
$ws.Cells.Item(5,5).Value = "new"

# Row 6 - Springvale (Shopping Centre)
$ws.Cells.Item(6,1).Value = "Springvale"
$ws.Cells.Item(6,2).Value = "Springvale Shopping Centre,  46-58 Buckingham Avenue"
$ws.Cells.Item(6,3).Value = "29/12/20, 11:00am-12:30pm"
$ws.Cells.Item(6,4).Value = "Case shopped"
$ws.Cells.Item(6,5).Value = "new"

# Column width adjustments to match new content (bestFit-style widths from target)
$ws.Columns.Item(1).ColumnWidth = 9.53125
$ws.Columns.Item(2).ColumnWidth = 56.9296875
$ws.Columns.Item(3).ColumnWidth = 24
$ws.Columns.Item(4).ColumnWidth = 32.3984375
